# Trade #188 closed at 2026-02-17 10:07:23 - unknown UNKNOWN +0.000%
#
# 1) New closed trade for "volatility_scorer" (Trade #188) appended to
#    "All Trades" and to the strategy's own "volatility_scorer" sheet.
# 2) New open trade for "MarketMaking" (Trade #189) appended to
#    "All Trades" and to the strategy's own "MarketMaking" sheet.
# 3) Summary / Strategy Status roll-up numbers updated to reflect the
#    newly closed trade.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.58   # Current Capital
$summary.Range("B4").Value = -0.42     # Total P&L $
$summary.Range("B5").Value = -0.04     # Total P&L %
$summary.Range("B6").Value = 188       # Total Trades
$summary.Range("B7").Value = 78        # Winning Trades
$summary.Range("B9").Value = 41.49     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - volatility_scorer row (row 12)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C12").Value = 99.45
$status.Range("D12").Value = 12
$status.Range("E12").Value = -0.55
$status.Range("F12").Value = -0.55
$status.Range("G12").Value = 41.67

# ---------------------------------------------------------------------
# All Trades sheet - append the two new trades (rows 189 and 190)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Row 189 = Trade #188 (volatility_scorer, CLOSED)
$allTrades.Range("A189").Value = 188
$allTrades.Range("B189").Value = "'2026-02-17"
$allTrades.Range("C189").Value = "10:07:16"
$allTrades.Range("D189").Value = "volatility_scorer"
$allTrades.Range("E189").Value = "NEUTRAL"
$allTrades.Range("F189").Value = 0.29
$allTrades.Range("G189").Value = 0.3
$allTrades.Range("H189").Value = "CLOSED"
$allTrades.Range("I189").Value = 3.4483
$allTrades.Range("J189").Value = 0.01
$allTrades.Range("K189").Value = 99.45
$allTrades.Range("L189").Value = 0
$allTrades.Range("M189").Value = 0
$allTrades.Range("N189").Value = 0.85
$allTrades.Range("O189").Value = "Low vol market (score: inf) - ideal for market making"
$allTrades.Range("P189").Value = "early_exit"
$allTrades.Range("Q189").Value = 0.17

# Row 190 = Trade #189 (MarketMaking, OPEN)
$allTrades.Range("A190").Value = 189
$allTrades.Range("B190").Value = "'2026-02-17"
$allTrades.Range("C190").Value = "10:07:16"
$allTrades.Range("D190").Value = "MarketMaking"
$allTrades.Range("E190").Value = "UP"
$allTrades.Range("F190").Value = 0.7
$allTrades.Range("G190").Value = "'"
$allTrades.Range("H190").Value = "OPEN"
$allTrades.Range("I190").Value = 0
$allTrades.Range("J190").Value = 0
$allTrades.Range("K190").Value = 100.1280687506789
$allTrades.Range("L190").Value = 0
$allTrades.Range("M190").Value = 0
$allTrades.Range("N190").Value = 0.6
$allTrades.Range("O190").Value = "Normal spread capture: 19600 bps"
$allTrades.Range("P190").Value = "'"
$allTrades.Range("Q190").Value = 0

# ---------------------------------------------------------------------
# volatility_scorer sheet - append its own copy of trade #188 (row 13)
# ---------------------------------------------------------------------
$volSheet = $wb.Worksheets.Item("volatility_scorer")
$volSheet.Range("A13").Value = 188
$volSheet.Range("B13").Value = "'2026-02-17"
$volSheet.Range("C13").Value = "10:07:16"
$volSheet.Range("D13").Value = "volatility_scorer"
$volSheet.Range("E13").Value = "NEUTRAL"
$volSheet.Range("F13").Value = 0.29
$volSheet.Range("G13").Value = 0.3
$volSheet.Range("H13").Value = "CLOSED"
$volSheet.Range("I13").Value = 3.4483
$volSheet.Range("J13").Value = 0.01
$volSheet.Range("K13").Value = 99.45
$volSheet.Range("L13").Value = 0
$volSheet.Range("M13").Value = 0
$volSheet.Range("N13").Value = 0.85
$volSheet.Range("O13").Value = "Low vol market (score: inf) - ideal for market making"
$volSheet.Range("P13").Value = "early_exit"
$volSheet.Range("Q13").Value = 0.17

# ---------------------------------------------------------------------
# MarketMaking sheet - append its own copy of trade #189 (row 178)
# ---------------------------------------------------------------------
$mmSheet = $wb.Worksheets.Item("MarketMaking")
$mmSheet.Range("A178").Value = 189
$mmSheet.Range("B178").Value = "'2026-02-17"
$mmSheet.Range("C178").Value = "10:07:16"
$mmSheet.Range("D178").Value = "MarketMaking"
$mmSheet.Range("E178").Value = "UP"
$mmSheet.Range("F178").Value = 0.7
$mmSheet.Range("G178").Value = "'"
$mmSheet.Range("H178").Value = "OPEN"
$mmSheet.Range("I178").Value = 0
$mmSheet.Range("J178").Value = 0
$mmSheet.Range("K178").Value = 100.1280687506789
$mmSheet.Range("L178").Value = 0
$mmSheet.Range("M178").Value = 0
$mmSheet.Range("N178").Value = 0.6
$mmSheet.Range("O178").Value = "Normal spread capture: 19600 bps"
$mmSheet.Range("P178").Value = "'"
$mmSheet.Range("Q178").Value = 0
